# Code review and refactoring of the CPACS/JSBSim comparison sheet:
#  - Correct/clarify the wording of the aileron deflection comment (B14)
#  - Add a new note about inverting deltaA in flight gear (B26)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = "in the CPACS defined until 3°, in JSBSim 40° (0.04 in CPACS (3°) and 0.9 (40°) in JSBSim) same in outer flap"
$ws.Range("B26").Value = "invertire deltaA in flight gear"

# Update the active cell selection to match the newly added row, as seen
# in the authoritative diff (selection moves from B25 to B26).
$ws.Range("B26").Select()
